$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9884.111000000001
$ws.Range("I69").Value = 9574.799999999999
$ws.Range("K69").Value = 28724.4
$ws.Range("M69").Value = -27850.4
$ws.Range("H72").Value = 9884.111000000001
$ws.Range("I72").Value = 9574.799999999999
$ws.Range("K72").Value = 86173.2
$ws.Range("M72").Value = -81805.2
$ws.Range("H86").Value = 4722.5713
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 4862.769
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 4862.769
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -7108.769
$ws.Range("H89").Value = 4722.5713
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 4862.769
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 24313.845
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -35545.845
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H98").Value = 2627.6875
$ws.Range("I98").Value = 2467.4285
$ws.Range("K98").Value = 2467.4285
$ws.Range("M98").Value = -969.4285
$ws.Range("H122").Value = 2627.6875
$ws.Range("I122").Value = 2467.4285
$ws.Range("K122").Value = 7402.2855
$ws.Range("M122").Value = -4952.2855
$ws.Range("H132").Value = 50403.957
$ws.Range("I132").Value = 51157.105
$ws.Range("K132").Value = 153471.315
$ws.Range("M132").Value = -150941.315
$ws.Range("H137").Value = 2021649.1
$ws.Range("I137").Value = 2316434.5
$ws.Range("K137").Value = 6949303.5
$ws.Range("M137").Value = -6946753.5
$ws.Range("H138").Value = 3256.1562
$ws.Range("I138").Value = 2033.7142
$ws.Range("J138").Value = 5589.909
$ws.Range("K138").Value = 6101.142599999999
$ws.Range("L138").Value = 16769.727
$ws.Range("M138").Value = -961.1425999999992
$ws.Range("N138").Value = -27049.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11118025
$ws.Range("I32").Value = 11909265
$ws.Range("K32").Value = 11909265
$ws.Range("M32").Value = -11908978
$ws.Range("H92").Value = 78875
$ws.Range("J92").Value = 78875
$ws.Range("L92").Value = 78875
$ws.Range("N92").Value = -83867
$ws.Range("H132").Value = 599119.4
$ws.Range("I132").Value = 647798.75
$ws.Range("K132").Value = 1943396.25
$ws.Range("M132").Value = -1940866.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1563.5862
$ws.Range("I20").Value = 1690
$ws.Range("K20").Value = 1690
$ws.Range("M20").Value = -1443
$ws.Range("H134").Value = 922683.0600000001
$ws.Range("I134").Value = 1080719.5
$ws.Range("K134").Value = 3242158.5
$ws.Range("M134").Value = -3239623.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4357
$ws.Range("I21").Value = 4357
$ws.Range("K21").Value = 4357
$ws.Range("M21").Value = -4122
$ws.Range("H31").Value = 4941
$ws.Range("I31").Value = 1770.091
$ws.Range("K31").Value = 1770.091
$ws.Range("M31").Value = -1475.091
$ws.Range("H34").Value = 4941
$ws.Range("I34").Value = 1770.091
$ws.Range("K34").Value = 1770.091
$ws.Range("M34").Value = -1568.091
$ws.Range("H86").Value = 3599.625
$ws.Range("I86").Value = 3633
$ws.Range("K86").Value = 3633
$ws.Range("M86").Value = -2510
$ws.Range("H89").Value = 3599.625
$ws.Range("I89").Value = 3633
$ws.Range("K89").Value = 18165
$ws.Range("M89").Value = -12549
$ws.Range("H134").Value = 1805.5264
$ws.Range("I134").Value = 1856.1666
$ws.Range("K134").Value = 5568.4998
$ws.Range("M134").Value = -3033.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1090
$ws.Range("J92").Value = 1553.5714
$ws.Range("L92").Value = 4660.7142
$ws.Range("N92").Value = -7156.7142
$ws.Range("H141").Value = 2868.0833
$ws.Range("I141").Value = 2868.0833
$ws.Range("K141").Value = 8604.249899999999
$ws.Range("M141").Value = -3424.249899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5365.5
$ws.Range("I70").Value = 5820.778
$ws.Range("K70").Value = 5820.778
$ws.Range("M70").Value = -5550.778
$ws.Range("H73").Value = 5365.5
$ws.Range("I73").Value = 5820.778
$ws.Range("K73").Value = 5820.778
$ws.Range("M73").Value = -4884.778
$ws.Range("H80").Value = 125472.36
$ws.Range("I80").Value = 191030.48
$ws.Range("J80").Value = 7467.7334
$ws.Range("K80").Value = 191030.48
$ws.Range("L80").Value = 7467.7334
$ws.Range("M80").Value = -190032.48
$ws.Range("N80").Value = -9463.733400000001
$ws.Range("H83").Value = 125472.36
$ws.Range("I83").Value = 191030.48
$ws.Range("J83").Value = 7467.7334
$ws.Range("K83").Value = 955152.4
$ws.Range("L83").Value = 37338.667
$ws.Range("M83").Value = -950160.4
$ws.Range("N83").Value = -47322.667
$ws.Range("H92").Value = 45046.332
$ws.Range("J92").Value = 45046.332
$ws.Range("L92").Value = 45046.332
$ws.Range("N92").Value = -48790.332
$ws.Range("H132").Value = 1508078.9
$ws.Range("I132").Value = 2009439.5
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 6028318.5
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -6025788.5
$ws.Range("N132").Value = -17051

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1489.4117
$ws.Range("I82").Value = 1267.4117
$ws.Range("J82").Value = 1711.4117
$ws.Range("K82").Value = 1267.4117
$ws.Range("L82").Value = 1711.4117
$ws.Range("M82").Value = -906.4117000000001
$ws.Range("N82").Value = -2433.4117
$ws.Range("H85").Value = 1489.4117
$ws.Range("I85").Value = 1267.4117
$ws.Range("J85").Value = 1711.4117
$ws.Range("K85").Value = 1267.4117
$ws.Range("L85").Value = 1711.4117
$ws.Range("M85").Value = -19.41170000000011
$ws.Range("N85").Value = -4207.411700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 96949.7
$ws.Range("J62").Value = 7721.8887
$ws.Range("L62").Value = 7721.8887
$ws.Range("N62").Value = -8969.8887
$ws.Range("H65").Value = 96949.7
$ws.Range("J65").Value = 7721.8887
$ws.Range("L65").Value = 38609.4435
$ws.Range("N65").Value = -44849.4435
$ws.Range("H68").Value = 71875.664
$ws.Range("J68").Value = 70196.8
$ws.Range("L68").Value = 70196.8
$ws.Range("N68").Value = -71818.8
$ws.Range("H71").Value = 71875.664
$ws.Range("J71").Value = 70196.8
$ws.Range("L71").Value = 210590.4
$ws.Range("N71").Value = -218702.4
$ws.Range("H136").Value = 9774941
$ws.Range("I136").Value = 11211565
$ws.Range("K136").Value = 33634695
$ws.Range("M136").Value = -33632145
